# Update candidate records per "Adding the extent report changes"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: George.Thompson1a973@gmail.com -> George.Thompson1a73@gmail.com, new phone
$ws.Range("G2").Value = "George.Thompson1a73@gmail.com"
$ws.Range("H2").Value = 8844885577

# Row 3: George.Thompson2@gmail.com8888888888 -> George.Thompson3@gmail.com8888888888, new phone
$ws.Range("G3").Value = "George.Thompson3@gmail.com8888888888"
$ws.Range("H3").Value = 6688774422

# Active cell / selection moves to G3 (last cell touched)
[void]$ws.Range("G3").Select()
